$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -190.03616681730207
$ws.Range("C2").Value = -386.15890840440017
$ws.Range("D2").Value = -536.53068034738487

$ws.Range("B3").Value = -329.81368043359379
$ws.Range("C3").Value = -304.74344595039679
$ws.Range("D3").Value = -564.59125077399528

$ws.Range("B4").Value = -400.906278166047
$ws.Range("C4").Value = -317.43595203653797
$ws.Range("D4").Value = -733.17960379220426

$ws.Range("B5").Value = -410.52063479944718
$ws.Range("C5").Value = -333.13969282022856
$ws.Range("D5").Value = -493.85412465835395

$ws.Range("B6").Value = -354.26011710232416
$ws.Range("C6").Value = -260.44063513311795
$ws.Range("D6").Value = -534.71259487172188

$ws.Range("B7").Value = -357.20350774086046
$ws.Range("C7").Value = -379.29475514103149
$ws.Range("D7").Value = -558.56169780240066

$ws.Range("B8").Value = -303.45368812069552
$ws.Range("C8").Value = -230.07951127841898
$ws.Range("D8").Value = -258.37987702984196

$ws.Range("B9").Value = -296.7614332654619
$ws.Range("C9").Value = -263.88575543447024
$ws.Range("D9").Value = -383.37207061921686

$ws.Range("A10").Value = 2023
$ws.Range("B10").Value = -308.15787901734899
$ws.Range("C10").Value = -279.86305368101188
$ws.Range("D10").Value = -517.20884682061887
